$wb = $excel.ActiveWorkbook

# --- Credentials sheet ---
$ws = $wb.Worksheets.Item("Credentials")
$ws.Range("A2").Value = "admin@testing.com"
$ws.Range("B2").Value = "truetesting@123"

# --- Email sheet ---
$ws = $wb.Worksheets.Item("Email")
$ws.Range("A2").Value = "newtester@tmail.com"

# --- ProductDetails sheet ---
$ws = $wb.Worksheets.Item("ProductDetails")
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2"

# --- AccountCreationData sheet ---
$ws = $wb.Worksheets.Item("AccountCreationData")
$ws.Range("A2").Value = "newTester1@kmail.com"
$ws.Range("A3").Value = "newTester2@kmail.com"
$ws.Range("A4").Value = "newTester3@kmail.com"
